# Split single runs of text into multiple <w:t> segments joined by
# manual line breaks (<w:br/>), matching the target diff.
#
# Using Find/Replace with "^l" (manual line break) in the replacement
# text causes Word to insert a real line-break character, which when
# saved back to OOXML serializes as </w:t><w:br/><w:t> inside the same
# run - exactly the structure required by the diff.

$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $range = $d.Content
    $result = $range.Find.Execute(
        $findText,   # FindText
        $true,       # MatchCase
        $false,      # MatchWholeWord
        $false,      # MatchWildcards
        $false,      # MatchSoundsLike
        $false,      # MatchAllWordForms
        $true,       # Forward
        1,           # Wrap (wdFindContinue)
        $false,      # Format
        $replaceText,# ReplaceWith
        2            # Replace (wdReplaceAll)
    )
    if (-not $result) {
        throw "Find.Execute failed for: $findText"
    }
}

# 1) "Programa resumido" - italic English description paragraph
Replace-Text "inglês.Real" "inglês.^lReal"

# 2) "Programa" - Portuguese bulleted paragraph
Replace-Text "inversa.•Limite" "inversa.^l•Limite"
Replace-Text "infinito.•Continuidade" "infinito.^l•Continuidade"
Replace-Text "intermediário.•Derivada" "intermediário.^l•Derivada"

# 3) "Programa" - English (italic) bulleted paragraph
Replace-Text "functions.•Limits" "functions.^l•Limits"
Replace-Text "infinite.•Continuity" "infinite.^l•Continuity"
Replace-Text "theorem.•Derivative" "theorem.^l•Derivative"
